$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: updated CV stats only (H2:K2) ---
$ws.Range("H2").Value = 0.5889630819184842
$ws.Range("I2").Value = 0.02037244586123459
$ws.Range("J2").Value = 0.4585093003401826
$ws.Range("K2").Value = 0.07788852758599678

# --- Row 3: pipeline now includes NamedFeatureSelector + class_weight=balanced ---
$A3 = @"
Pipeline(steps=[('scaler', RobustScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f91ca0dc4f0>),
                ('model',
                 SVC(C=1, class_weight='balanced', kernel='sigmoid',
                     random_state=42))])
"@
$ws.Range("A3").Value = $A3
$ws.Range("B3").Value = 0.7565567765567764
$C3 = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f91c9d73580>, 'scaler': RobustScaler(), 'model__kernel': 'sigmoid', 'model__class_weight': 'balanced', 'model__C': 1}
"@
$ws.Range("C3").Value = $C3
$ws.Range("D3").Value = 0.75
$ws.Range("F3").Value = "[1 1 1 1 1 0 1 1 1 0 1 0]"
$ws.Range("H3").Value = 0.6850921714225813
$ws.Range("I3").Value = 0.01953896842368958
$ws.Range("J3").Value = 0.5841354729420906
$ws.Range("K3").Value = 0.06679899670194876

# --- Row 4: CV stats only (H4:K4) ---
$ws.Range("H4").Value = 0.6281890442110398
$ws.Range("I4").Value = 0.01963124591775743
$ws.Range("J4").Value = 0.541458159569189
$ws.Range("K4").Value = 0.07022563240921781

# --- Row 5: pipeline now uses RobustScaler + NamedFeatureSelector, no class_weight ---
$A5 = @"
Pipeline(steps=[('scaler', RobustScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f91c9d8d820>),
                ('model', SVC(C=5, kernel='sigmoid', random_state=42))])
"@
$ws.Range("A5").Value = $A5
$ws.Range("B5").Value = 0.7954761904761904
$C5 = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f91ca1954c0>, 'scaler': RobustScaler(), 'model__kernel': 'sigmoid', 'model__class_weight': None, 'model__C': 5}
"@
$ws.Range("C5").Value = $C5
$ws.Range("D5").Value = 0.7692307692307693
$ws.Range("F5").Value = "[0 1 1 0 1 1 0 0 1 1 1 1]"
$ws.Range("H5").Value = 0.7510527510511276
$ws.Range("I5").Value = 0.01922681041905444
$ws.Range("J5").Value = 0.7083854708200298
$ws.Range("K5").Value = 0.07029864901140509

